# Adds the "bulk add" transactions captured on 2024-11-17 (evening) across
# the All_SANs master list plus the per-site Items/Timestamps sheets, so
# staff can enter a batch of new SANs faster without re-keying totals.

$wb = $excel.ActiveWorkbook

$RDQ = [char]0x201D   # right double quotation mark, used in "Monitor 24<RDQ>"
$NBSP = [char]0x00A0  # trailing non-breaking space used after Monitor 24<RDQ>
$Monitor24 = "Monitor 24" + $RDQ + $NBSP

# ---------------------------------------------------------------------
# All_SANs: append the 12 newly-logged serial numbers
# ---------------------------------------------------------------------
$allSans = $wb.Worksheets.Item("All_SANs")

$newSans = @(
  @("SAN444444", "Laptop 840 G10", "2024-11-17 18:44:41", "BR"),
  @("SAN343435", "Laptop 840 G10", "2024-11-17 18:49:44", "BR"),
  @("SAN555444", "Laptop 840 G10", "2024-11-17 18:49:46", "BR"),
  @("SAN45453",  "Laptop 840 G10", "2024-11-17 18:49:54", "BR"),
  @("SAN456753", "Laptop x360 G8", "2024-11-17 18:52:31", "4.2"),
  @("SAN111567", "Laptop x360 G8", "2024-11-17 18:52:32", "4.2"),
  @("SAN494946", "Laptop x360 G8", "2024-11-17 18:52:37", "4.2"),
  @("SAN166544", "Laptop x360 G8", "2024-11-17 18:52:40", "4.2"),
  @("SAN493574", "Laptop x360 G8", "2024-11-17 18:56:13", "4.2"),
  @("SAN165478", "Laptop x360 G8", "2024-11-17 18:56:16", "4.2"),
  @("SAN254687", "Laptop x360 G8", "2024-11-17 18:56:18", "4.2"),
  @("SAN264578", "Laptop x360 G8", "2024-11-17 18:56:19", "4.2")
)

$startRow = 132
$endRow = $startRow + $newSans.Length - 1

# Force column D to store as text — the "4.2" site code would otherwise be
# auto-coerced into a number by the normal Value setter.
$allSans.Range("D" + $startRow + ":D" + $endRow).NumberFormat = "@"

for ($i = 0; $i -lt $newSans.Length; $i++) {
  $row = $startRow + $i
  $rec = $newSans[$i]
  $allSans.Cells.Item($row, 1).Value = $rec[0]
  $allSans.Cells.Item($row, 2).Value = $rec[1]
  $allSans.Cells.Item($row, 3).Value = $rec[2]
  $allSans.Cells.Item($row, 4).Value = $rec[3]
}
# New rows should stay unstyled (matching the rest of the recent entries),
# not inherit the centred column style used further up the sheet, and not
# keep the "@" text format applied above.
$allSans.Range("A132:D143").Style = "Normal"

# ---------------------------------------------------------------------
# 4.2_Items: Laptop x360 G8 totals
# ---------------------------------------------------------------------
$items42 = $wb.Worksheets.Item("4.2_Items")
$items42.Cells.Item(10, 2).Value = 18
$items42.Cells.Item(10, 3).Value = 22

# ---------------------------------------------------------------------
# 4.2_Timestamps: log each add + the rollup "add 4" entries
# ---------------------------------------------------------------------
$ts42 = $wb.Worksheets.Item("4.2_Timestamps")
$ts42.Cells.Item(32, 4).ClearContents()

$ts42.Cells.Item(33, 1).Value = "2024-11-17 18:52:31"
$ts42.Cells.Item(33, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(33, 3).Value = "add"
$ts42.Cells.Item(33, 4).Value = "SAN456753"

$ts42.Cells.Item(34, 1).Value = "2024-11-17 18:52:32"
$ts42.Cells.Item(34, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(34, 3).Value = "add"
$ts42.Cells.Item(34, 4).Value = "SAN111567"

$ts42.Cells.Item(35, 1).Value = "2024-11-17 18:52:37"
$ts42.Cells.Item(35, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(35, 3).Value = "add"
$ts42.Cells.Item(35, 4).Value = "SAN494946"

$ts42.Cells.Item(36, 1).Value = "2024-11-17 18:52:40"
$ts42.Cells.Item(36, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(36, 3).Value = "add"
$ts42.Cells.Item(36, 4).Value = "SAN166544"

$ts42.Cells.Item(37, 1).Value = "2024-11-17 18:52:40"
$ts42.Cells.Item(37, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(37, 3).Value = "add 4"

$ts42.Cells.Item(38, 1).Value = "2024-11-17 18:56:13"
$ts42.Cells.Item(38, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(38, 3).Value = "add"
$ts42.Cells.Item(38, 4).Value = "SAN493574"

$ts42.Cells.Item(39, 1).Value = "2024-11-17 18:56:16"
$ts42.Cells.Item(39, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(39, 3).Value = "add"
$ts42.Cells.Item(39, 4).Value = "SAN165478"

$ts42.Cells.Item(40, 1).Value = "2024-11-17 18:56:18"
$ts42.Cells.Item(40, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(40, 3).Value = "add"
$ts42.Cells.Item(40, 4).Value = "SAN254687"

$ts42.Cells.Item(41, 1).Value = "2024-11-17 18:56:19"
$ts42.Cells.Item(41, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(41, 3).Value = "add"
$ts42.Cells.Item(41, 4).Value = "SAN264578"

$ts42.Cells.Item(42, 1).Value = "2024-11-17 18:56:20"
$ts42.Cells.Item(42, 2).Value = "Laptop x360 G8"
$ts42.Cells.Item(42, 3).Value = "add 4"

# ---------------------------------------------------------------------
# BR_Items: Laptop 840 G10 totals
# ---------------------------------------------------------------------
$itemsBR = $wb.Worksheets.Item("BR_Items")
$itemsBR.Cells.Item(6, 2).Value = 1
$itemsBR.Cells.Item(6, 3).Value = 4

# ---------------------------------------------------------------------
# BR_Timestamps: log each add + the rollup entries
# ---------------------------------------------------------------------
$tsBR = $wb.Worksheets.Item("BR_Timestamps")

$tsBR.Cells.Item(22, 1).Value = "2024-11-17 18:44:41"
$tsBR.Cells.Item(22, 2).Value = "Laptop 840 G10"
$tsBR.Cells.Item(22, 3).Value = "add"
$tsBR.Cells.Item(22, 4).Value = "SAN444444"

$tsBR.Cells.Item(23, 1).Value = "2024-11-17 18:44:41"
$tsBR.Cells.Item(23, 2).Value = "Laptop 840 G10"
$tsBR.Cells.Item(23, 3).Value = "add 1"

$tsBR.Cells.Item(24, 1).Value = "2024-11-17 18:49:44"
$tsBR.Cells.Item(24, 2).Value = "Laptop 840 G10"
$tsBR.Cells.Item(24, 3).Value = "add"
$tsBR.Cells.Item(24, 4).Value = "SAN343435"

$tsBR.Cells.Item(25, 1).Value = "2024-11-17 18:49:46"
$tsBR.Cells.Item(25, 2).Value = "Laptop 840 G10"
$tsBR.Cells.Item(25, 3).Value = "add"
$tsBR.Cells.Item(25, 4).Value = "SAN555444"

$tsBR.Cells.Item(26, 1).Value = "2024-11-17 18:49:54"
$tsBR.Cells.Item(26, 2).Value = "Laptop 840 G10"
$tsBR.Cells.Item(26, 3).Value = "add"
$tsBR.Cells.Item(26, 4).Value = "SAN45453"

$tsBR.Cells.Item(27, 1).Value = "2024-11-17 18:49:54"
$tsBR.Cells.Item(27, 2).Value = "Laptop 840 G10"
$tsBR.Cells.Item(27, 3).Value = "add 3"

# ---------------------------------------------------------------------
# L17_Items: Laptop 840 G6 + Monitor 24" totals
# ---------------------------------------------------------------------
$itemsL17 = $wb.Worksheets.Item("L17_Items")
$itemsL17.Cells.Item(2, 2).Value = 15
$itemsL17.Cells.Item(2, 3).Value = 39
$itemsL17.Cells.Item(3, 2).Value = 3
$itemsL17.Cells.Item(3, 3).Value = 27

# ---------------------------------------------------------------------
# L17_Timestamps: log the two rollup adds
# ---------------------------------------------------------------------
$tsL17 = $wb.Worksheets.Item("L17_Timestamps")

$tsL17.Cells.Item(6, 1).Value = "2024-11-17 18:43:11"
$tsL17.Cells.Item(6, 2).Value = $Monitor24
$tsL17.Cells.Item(6, 3).Value = "add 24"

$tsL17.Cells.Item(7, 1).Value = "2024-11-17 18:43:14"
$tsL17.Cells.Item(7, 2).Value = "Laptop 840 G6"
$tsL17.Cells.Item(7, 3).Value = "add 24"

# ---------------------------------------------------------------------
# Darwin_Items: Monitor 24" + Wired Mouse totals
# ---------------------------------------------------------------------
$itemsDarwin = $wb.Worksheets.Item("Darwin_Items")
$itemsDarwin.Cells.Item(11, 2).Value = 30
$itemsDarwin.Cells.Item(11, 3).Value = 54
$itemsDarwin.Cells.Item(16, 2).Value = 0
$itemsDarwin.Cells.Item(16, 3).Value = 24

# ---------------------------------------------------------------------
# Darwin_Timestamps: log the two rollup adds
# ---------------------------------------------------------------------
$tsDarwin = $wb.Worksheets.Item("Darwin_Timestamps")
$tsDarwin.Cells.Item(28, 4).ClearContents()

$tsDarwin.Cells.Item(29, 1).Value = "2024-11-17 18:42:55"
$tsDarwin.Cells.Item(29, 2).Value = $Monitor24
$tsDarwin.Cells.Item(29, 3).Value = "add 24"

$tsDarwin.Cells.Item(30, 1).Value = "2024-11-17 18:43:00"
$tsDarwin.Cells.Item(30, 2).Value = "Wired Mouse"
$tsDarwin.Cells.Item(30, 3).Value = "add 24"
